# Update the "Uso_solo_simplificado" category (column C) to "Transporte"
# for a specific set of rows, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(14,15,16,17,38,39,40,41,54,55,56,57,66,67,68,69,74,75,76,77,82,83,84,85,102,103,104,105,110,111,112,113,126,127,128,129,142,143,144,145,146,147,148,149,150,151,152,153,162,163,164,165,174,175,176,177,178,179,180,181)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Transporte"
}
